# Update cryptocurrency price/volume data per upstream refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '42.316.86'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  -0.18%  '

# Row 3
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.246.33'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  -1.03%  '

# Row 4
$ws.Range('E4').Value = '  -0.03%  '

# Row 5
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '247.79'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -1.03%  '

# Row 6
$ws.Range('E6').Value = '  +0.05%  '

# Row 7
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '76.19'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  +6.62%  '

# Row 8
$ws.Range('E8').Value = '  +0.17%  '

# Row 9
$ws.Range('E9').Value = '  -1.39%  '

# Row 10
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '40.06'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +4.17%  '

# Row 11
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0952'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -1.24%  '

# Row 12
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '7.23'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -0.74%  '

# Row 13
$ws.Range('E13').Value = '  -0.85%  '

# Row 14
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '2.583.07'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -1.12%  '

# Row 15
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '14.89'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -0.12%  '

# Row 16
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.861'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -1.96%  '

# Row 17
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '2.248.42'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -0.42%  '

# Row 18
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '42.220.45'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -0.36%  '

# Row 19
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.0₃0980'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -0.99%  '

# Row 20
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '6.15'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -1.77%  '

# Row 21
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '71.48'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -1.77%  '

# Row 22
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '231.42'
$ws.Range('D22').ClearFormats()

# Row 23
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '2.17'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -4.24%  '

# Row 24
$ws.Range('E24').Value = '  -0.10%  '

# Row 25
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '3.72'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -5.11%  '

# Row 26
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '11.16'
$ws.Range('D26').ClearFormats()

# Row 27
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '2.32'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -4.71%  '

# Row 28
$ws.Range('E28').Value = '  -0.74%  '

# Row 29
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '6.89'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +7.35%  '

# Row 30
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '168.00'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +0.90%  '

# Row 31
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '20.54'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -1.86%  '

# Row 32
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.0853'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +7.03%  '

# Row 33
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '30.95'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -1.83%  '

# Row 34
$ws.Range('E34').Value = '  -5.59%  '

# Row 35
$ws.Range('E35').Value = '  +0.21%  '

# Row 36
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '4.48'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -6.45%  '

# Row 37
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '4.73'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +0.05%  '

# Row 38
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.0298'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -3.10%  '

# Row 39
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '12.98'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -4.56%  '

# Row 40
$ws.Range('E40').Value = '  -4.18%  '

# Row 41
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '5.92'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -1.05%  '

# Row 42
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '118.02'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +23.30%  '

# Row 43
$ws.Range('E43').Value = '  -2.94%  '

# Row 44
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '60.11'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -2.26%  '

# Row 45
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '8.74'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -5.49%  '

# Row 46
$ws.Range('E46').Value = '  -2.02%  '

# Row 47
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.996'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -0.34%  '

# Row 48
$ws.Range('E48').Value = '  -3.44%  '

# Row 49
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.17'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -1.00%  '

# Row 50
$ws.Range('B50').Value = 'FTXToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '4.16'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -14.95%  '

# Row 51
$ws.Range('B51').Value = 'SynthetixNetwork'
$ws.Range('C51').Value = 'https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '4.11'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -2.35%  '

Write-Output "Applied cryptos update"
